$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$qCases = @'
 MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis) 
 MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Female'] and demo.neutered_indicator IN ['Yes']
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
Order by c.case_id LIMIT 100        
'@
$qSamples = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Female'] and demo.neutered_indicator IN ['Yes']
 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$qFiles = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Female'] and demo.neutered_indicator IN ['Yes']
WITH DISTINCT f, parent, c, demo, diag, s
OPTIONAL MATCH (f)-[*]->(samp:sample)
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN coalesce(f.file_name, '') AS `File Name`, 
 coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`, 
      CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
   coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis 
        Order By f.file_name LIMIT 100
'@
$qStudyFiles = @'
MATCH (f:file)-->(s:study)
MATCH (s)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (sf:file)-->(s)
MATCH (s)<--(c)
MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Female'] and demo.neutered_indicator IN ['Yes']
WITH DISTINCT f,  s, c
WITH
        f, c,  s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c,  s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c,   s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$qStat = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (f:file)-[*]->(samp:sample)-->(c)
MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp,demo, c, s, p, diag
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in [ 'Bluetick Hound', 'Norfolk Terrier', 'Scottish Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Female'] and demo.neutered_indicator IN ['Yes']
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Row 2: CasesTab
$ws.Range("B2").Value = $qCases
$ws.Range("C2").Value = $qStat

# Row 3: SamplesTab
$ws.Range("B3").Value = $qSamples
$ws.Range("C3").Value = $qStat

# Row 4: FilesTab
$ws.Range("B4").Value = $qFiles
$ws.Range("C4").Value = $qStat

# Row 5: StudyFilesTab (new row)
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5:C5").WrapText = $true
$ws.Range("B5").Value = $qStudyFiles
$ws.Range("C5").Value = $qStat
$ws.Range("D5").Value = "TC21_Canine_StudyUBC02-Breed_Sex_NeuteredStatus_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC21_Canine_StudyUBC02-Breed_Sex_NeuteredStatus_WebData.xlsx"

# Update selection to match the final saved state (C5 selected, no frozen top-left override)
$ws.Range("C5").Select() | Out-Null
